$d = $word.ActiveDocument

# "dependiendo de en que posición" -> "dependiendo de en qué posición"
# (adds the missing accent mark on "que" -> "qué" in that one spot)
$d.Content.Find.Execute(
    "dependiendo de en que posición",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "dependiendo de en qué posición",
    2
)
